$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet
$ws.Name = "August"

# Header row
$ws.Range("A1").Value = "Event ID"
$ws.Range("B1").Value = "Type"
$ws.Range("C1").Value = "Title"
$ws.Range("D1").Value = "Start Date"
$ws.Range("E1").Value = "End Date"
$ws.Range("F1").Value = "Rooms"
$ws.Range("G1").Value = "Booked By"

# Row 2
$ws.Range("A2").Value = "-LKKaZGZAHjHCwBCkhPW"
$ws.Range("B2").Value = "External"
$ws.Range("C2").Value = "Web Dev Workshop"
$ws.Range("D2").Value = "Tuesday, 21st August 2018"
$ws.Range("E2").Value = "Sunday, 26th August 2018"
$ws.Range("F2").Value = "AB-5-201, AB-5-202, AB-5-203, AB-5-204"
$ws.Range("G2").Value = "Bhawesh Bhansali"

# Row 3
$ws.Range("A3").Value = "-LKKe4SAl7bWFfZVPEQz"
$ws.Range("B3").Value = "External"
$ws.Range("C3").Value = "Rotaract Work"
$ws.Range("D3").Value = "Tuesday, 21st August 2018"
$ws.Range("E3").Value = "Thursday, 23rd August 2018"
$ws.Range("F3").Value = "NLH-201, NLH-202"
$ws.Range("G3").Value = "Daksh"

# Row 4
$ws.Range("A4").Value = "dummcook46891"
$ws.Range("B4").Value = "Internal"
$ws.Range("C4").Value = "Cooking"
$ws.Range("D4").Value = "Tuesday, 21st August 2018"
$ws.Range("E4").Value = "Thursday, 23rd August 2018"
$ws.Range("F4").Value = "AB-5-310, AB-5-311"
$ws.Range("G4").Value = "Bhawesh"

# Row 5
$ws.Range("A5").Value = "dummlits594030"
$ws.Range("B5").Value = "Internal"
$ws.Range("C5").Value = "litstock"
$ws.Range("D5").Value = "Thursday, 23rd August 2018"
$ws.Range("E5").Value = "Saturday, 25th August 2018"
$ws.Range("F5").Value = "NLH-303, NLH-305"
$ws.Range("G5").Value = "vibhuti"

# Row 6
$ws.Range("A6").Value = "dummplan848023"
$ws.Range("B6").Value = "External"
$ws.Range("C6").Value = "Planning Workshop"
$ws.Range("D6").Value = "Thursday, 23rd August 2018"
$ws.Range("E6").Value = "Friday, 24th August 2018"
$ws.Range("F6").Value = "AB-5-201, AB-5-202"
$ws.Range("G6").Value = "Bhawesh"

# Column widths
$ws.Columns.Item(1).ColumnWidth = 25
$ws.Columns.Item(2).ColumnWidth = 10
$ws.Range("C1:F1").ColumnWidth = 25
$ws.Columns.Item(7).ColumnWidth = 15
